# Scheduled-runner update: refresh Universalis market-price-derived leve
# profit figures (currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Only columns H:N (price/profit
# columns) change per row; item/leve identity columns A:G are untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2477.1428
$ws.Range("I100").Value = 2640.25
$ws.Range("J100").Value = 1498.5
$ws.Range("K100").Value = 2640.25
$ws.Range("L100").Value = 1498.5
$ws.Range("M100").Value = -2099.25
$ws.Range("N100").Value = -2580.5

$ws.Range("H138").Value = 3410.4
$ws.Range("J138").Value = 5199.6665
$ws.Range("L138").Value = 15598.9995
$ws.Range("N138").Value = -25878.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""

$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H32").Value = 5267.077
$ws.Range("I32").Value = 3187.8572
$ws.Range("K32").Value = 3187.8572
$ws.Range("M32").Value = -2900.8572

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = ""
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = 0

$ws.Range("H63").Value = 4397.25
$ws.Range("I63").Value = 2655.4
$ws.Range("J63").Value = 5641.4287
$ws.Range("K63").Value = 2655.4
$ws.Range("L63").Value = 5641.4287
$ws.Range("M63").Value = -1969.4
$ws.Range("N63").Value = -7013.4287

$ws.Range("H66").Value = 4397.25
$ws.Range("I66").Value = 2655.4
$ws.Range("J66").Value = 5641.4287
$ws.Range("K66").Value = 13277
$ws.Range("L66").Value = 28207.1435
$ws.Range("M66").Value = -9845
$ws.Range("N66").Value = -35071.14350000001

$ws.Range("H97").Value = 1356.25
$ws.Range("I97").Value = 939.2308
$ws.Range("J97").Value = 3163.3333
$ws.Range("K97").Value = 939.2308
$ws.Range("L97").Value = 3163.3333
$ws.Range("M97").Value = -443.2308
$ws.Range("N97").Value = -4155.3333

$ws.Range("H110").Value = 1053.2
$ws.Range("I110").Value = 1053.2
$ws.Range("K110").Value = 1053.2
$ws.Range("M110").Value = 991.8

$ws.Range("H122").Value = 3027.6365
$ws.Range("J122").Value = 4661.6665
$ws.Range("L122").Value = 13984.9995
$ws.Range("N122").Value = -18884.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = ""

$ws.Range("H82").Value = 4932.75
$ws.Range("I82").Value = 4932.75
$ws.Range("K82").Value = 4932.75
$ws.Range("M82").Value = -4549.75

$ws.Range("H85").Value = 4932.75
$ws.Range("I85").Value = 4932.75
$ws.Range("K85").Value = 4932.75
$ws.Range("M85").Value = -3606.75

$ws.Range("H94").Value = 4140.6875
$ws.Range("I94").Value = 3750.0908
$ws.Range("K94").Value = 3750.0908
$ws.Range("M94").Value = -3299.0908

$ws.Range("H134").Value = 2972.4614
$ws.Range("I134").Value = 3158.818
$ws.Range("J134").Value = 1947.5
$ws.Range("K134").Value = 9476.454000000002
$ws.Range("L134").Value = 5842.5
$ws.Range("M134").Value = -6941.454000000002
$ws.Range("N134").Value = -10912.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2189.3333
$ws.Range("I22").Value = 546.3333
$ws.Range("J22").Value = 3832.3333
$ws.Range("K22").Value = 546.3333
$ws.Range("L22").Value = 3832.3333
$ws.Range("M22").Value = -196.3333
$ws.Range("N22").Value = -4532.3333

$ws.Range("H31").Value = 2807.111
$ws.Range("I31").Value = 3035
$ws.Range("K31").Value = 3035
$ws.Range("M31").Value = -2740

$ws.Range("H34").Value = 2807.111
$ws.Range("I34").Value = 3035
$ws.Range("K34").Value = 3035
$ws.Range("M34").Value = -2833

$ws.Range("H58").Value = 2595.4
$ws.Range("I58").Value = 2050
$ws.Range("K58").Value = 2050
$ws.Range("M58").Value = -1847

$ws.Range("H107").Value = 360
$ws.Range("I107").Value = 343
$ws.Range("J107").Value = 368.5
$ws.Range("K107").Value = 343
$ws.Range("L107").Value = 368.5
$ws.Range("M107").Value = 1577
$ws.Range("N107").Value = -4208.5

$ws.Range("H122").Value = 1284.3334
$ws.Range("I122").Value = 1209.6
$ws.Range("J122").Value = 1377.75
$ws.Range("K122").Value = 3628.8
$ws.Range("L122").Value = 4133.25
$ws.Range("M122").Value = -1178.8
$ws.Range("N122").Value = -9033.25

$ws.Range("H134").Value = 1831.3334
$ws.Range("I134").Value = 1831.3334
$ws.Range("K134").Value = 5494.0002
$ws.Range("M134").Value = -2959.0002

$ws.Range("H136").Value = 2595.4
$ws.Range("I136").Value = 2050
$ws.Range("K136").Value = 6150
$ws.Range("M136").Value = -3600

$ws.Range("H141").Value = 399979
$ws.Range("J141").Value = 399979
$ws.Range("L141").Value = 399979
$ws.Range("N141").Value = -410339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 53.875
$ws.Range("I10").Value = 66.75
$ws.Range("J10").Value = 41
$ws.Range("K10").Value = 200.25
$ws.Range("L10").Value = 123
$ws.Range("M10").Value = -61.25
$ws.Range("N10").Value = -401

$ws.Range("H15").Value = 174.75
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 174.75
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = 524.25
$ws.Range("N15").Value = -804.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 575.6
$ws.Range("I2").Value = 2436.5
$ws.Range("J2").Value = 190.58621
$ws.Range("K2").Value = 2436.5
$ws.Range("L2").Value = 190.58621
$ws.Range("M2").Value = -2323.5
$ws.Range("N2").Value = -416.58621

$ws.Range("H68").Value = 45000
$ws.Range("I68").Value = 45000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 45000
$ws.Range("L68").Value = ""
$ws.Range("M68").Value = -44189
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 45000
$ws.Range("I71").Value = 45000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 135000
$ws.Range("L71").Value = ""
$ws.Range("M71").Value = -130944
$ws.Range("N71").Value = 0

$ws.Range("H97").Value = 1124.7778
$ws.Range("I97").Value = 738
$ws.Range("J97").Value = 1898.3334
$ws.Range("K97").Value = 738
$ws.Range("L97").Value = 1898.3334
$ws.Range("M97").Value = -242
$ws.Range("N97").Value = -2890.3334

$ws.Range("H107").Value = 322.5
$ws.Range("I107").Value = 270.3846
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 270.3846
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1649.6154
$ws.Range("N107").Value = -4840

$ws.Range("H113").Value = 2195
$ws.Range("I113").Value = 2195
$ws.Range("K113").Value = 2195
$ws.Range("M113").Value = -25

$ws.Range("H122").Value = 1089
$ws.Range("I122").Value = 1114.2727
$ws.Range("K122").Value = 3342.8181
$ws.Range("M122").Value = -892.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3421.9033
$ws.Range("I40").Value = 3231.4285
$ws.Range("K40").Value = 3231.4285
$ws.Range("M40").Value = -3095.4285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5533
$ws.Range("I62").Value = 4799.5
$ws.Range("K62").Value = 4799.5
$ws.Range("M62").Value = -4175.5

$ws.Range("H65").Value = 5533
$ws.Range("I65").Value = 4799.5
$ws.Range("K65").Value = 23997.5
$ws.Range("M65").Value = -20877.5
